$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-07-18 Thursday" "2024-07-19 Friday"

Replace-Text "653÷3=217, 2" "361÷5=72, 1"
Replace-Text "902÷6=150, 2" "391÷4=97, 3"
Replace-Text "796÷9=88, 4" "791÷9=87, 8"
Replace-Text "595÷7=85, 0" "261÷6=43, 3"
Replace-Text "899÷8=112, 3" "408÷7=58, 2"

Replace-Text "103÷8=12, 7" "263÷3=87, 2"
Replace-Text "983÷3=327, 2" "344÷9=38, 2"
Replace-Text "972÷7=138, 6" "956÷8=119, 4"
Replace-Text "957÷6=159, 3" "183÷2=91, 1"
Replace-Text "893÷8=111, 5" "697÷5=139, 2"

Replace-Text "106÷5=21, 1" "660÷4=165, 0"
Replace-Text "628÷5=125, 3" "807÷6=134, 3"
Replace-Text "863÷4=215, 3" "160÷5=32, 0"
Replace-Text "665÷6=110, 5" "253÷7=36, 1"
Replace-Text "597÷4=149, 1" "175÷9=19, 4"

Replace-Text "763÷7=109, 0" "315÷4=78, 3"
Replace-Text "317÷3=105, 2" "502÷8=62, 6"
Replace-Text "418÷7=59, 5" "923÷6=153, 5"
Replace-Text "585÷7=83, 4" "519÷5=103, 4"
Replace-Text "989÷8=123, 5" "921÷8=115, 1"

Replace-Text "187÷3=62, 1" "584÷6=97, 2"
Replace-Text "999÷4=249, 3" "454÷2=227, 0"
Replace-Text "915÷5=183, 0" "121÷9=13, 4"
Replace-Text "453÷5=90, 3" "120÷8=15, 0"
Replace-Text "497÷5=99, 2" "464÷5=92, 4"
